$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely, shifting the remaining columns (B:F -> A:E) to the left.
$ws.Columns("A").Delete()
